$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple "Selger" (vendor) value corrections: Digikey/DIgikey -> Farnell ---
$ws.Range("G2").Value  = "Farnell"
$ws.Range("G5").Value  = "Farnell"
$ws.Range("G6").Value  = "Farnell"
$ws.Range("G7").Value  = "Farnell"
$ws.Range("G9").Value  = "Farnell"
$ws.Range("G10").Value = "Farnell"
$ws.Range("G11").Value = "Farnell"

# --- Other part-number / text corrections ---
$ws.Range("F7").Value  = "5988A10107F"
$ws.Range("B8").Value  = "WR-MPC4 4.2mm Male Single Row Angled Header with Mounting Flanges for Screw-in Retention , 2p (male connnector)"
$ws.Range("F10").Value = "RN73H2ATTD1001B25"
$ws.Range("C11").Value = "R3"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = "MCWR08X0000FTL"
$ws.Range("B12").Value = "Pre crimped wires"
$ws.Range("F12").Value = "we: 649500116015"
$ws.Range("G12").Value = "Wurth electonics"

# B8 and C11 originally used the "quote prefix" style (matches A8/A11); plain
# Value assignment on these resets them to the unprefixed border style, so
# restore it from an untouched sibling cell in the same row.
$ws.Range("A8").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("A11").Copy()
$ws.Range("C11").PasteSpecial(-4122)

# --- New comment cell next to the TVS diode row ---
$ws.Range("H4").Value  = "<- Finner kun på farnell"

# --- New row 14: solar panel line item ---
$ws.Range("B14").Value = "Solcelle"
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = "https://voltaicsystems.com/3-5-watt-panel/"
$ws.Range("G14").Value = "Voltaic"

# --- Formatting: F7, G7, F10, F11 get a smaller plain (unfilled) black Arial style ---
$fmtSrc = $ws.Range("F7")
$fmtSrc.Font.Name = "Arial"
$fmtSrc.Font.Size = 9
$fmtSrc.Font.Color = 0
$fmtSrc.Interior.Pattern = -4142
$fmtSrc.HorizontalAlignment = -4142
$fmtSrc.VerticalAlignment = -4108
$fmtSrc.WrapText = $true

$fmtSrc.Copy()
$ws.Range("G7").PasteSpecial(-4122)
$fmtSrc.Copy()
$ws.Range("F10").PasteSpecial(-4122)
$fmtSrc.Copy()
$ws.Range("F11").PasteSpecial(-4122)

# --- Formatting: F12 reverts to the plain "part number" style used elsewhere (e.g. F3/F13) ---
$ws.Range("F3").Copy()
$ws.Range("F12").PasteSpecial(-4122)

# --- Formatting: new row 14 cells adopt styles matching the existing rows above them ---
$ws.Range("E2").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("G2").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("G2").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("G2").Copy()
$ws.Range("G14").PasteSpecial(-4122)

Write-Host "Edit applied"
